$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 39672.45
$ws.Range("I129").Value = 313.57144
$ws.Range("J129").Value = 60865.69
$ws.Range("K129").Value = 940.71432
$ws.Range("L129").Value = 182597.07
$ws.Range("M129").Value = 4059.28568
$ws.Range("N129").Value = -192597.07

$ws.Range("H132").Value = 19824578
$ws.Range("I132").Value = 30002552
$ws.Range("J132").Value = 4313.6313
$ws.Range("K132").Value = 90007656
$ws.Range("L132").Value = 12940.8939
$ws.Range("M132").Value = -90005126
$ws.Range("N132").Value = -18000.8939

$ws.Range("H135").Value = 741921.3
$ws.Range("I135").Value = 1129.3684
$ws.Range("J135").Value = 2021471.1
$ws.Range("K135").Value = 10164.3156
$ws.Range("L135").Value = 18193239.9
$ws.Range("M135").Value = -7629.3156
$ws.Range("N135").Value = -18198309.9

$ws.Range("H138").Value = 1860.7847
$ws.Range("I138").Value = 1226.4642
$ws.Range("J138").Value = 2340.8108
$ws.Range("K138").Value = 3679.3926
$ws.Range("L138").Value = 7022.432400000001
$ws.Range("M138").Value = 1460.6074
$ws.Range("N138").Value = -17302.4324

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 63266.047
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 63266.047
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 63266.047
$ws.Range("M32").ClearContents()
$ws.Range("N32").Value = -63840.047

$ws.Range("H61").Value = 5875.4346
$ws.Range("I61").Value = 7748.923
$ws.Range("J61").Value = 3439.9
$ws.Range("K61").Value = 7748.923
$ws.Range("L61").Value = 3439.9
$ws.Range("M61").Value = -7536.923
$ws.Range("N61").Value = -3863.9

$ws.Range("H82").Value = 40000
$ws.Range("J82").Value = 40000
$ws.Range("L82").Value = 40000
$ws.Range("N82").Value = -40722

$ws.Range("H85").Value = 40000
$ws.Range("J85").Value = 40000
$ws.Range("L85").Value = 40000
$ws.Range("N85").Value = -42496

$ws.Range("H136").Value = 5875.4346
$ws.Range("I136").Value = 7748.923
$ws.Range("J136").Value = 3439.9
$ws.Range("K136").Value = 23246.769
$ws.Range("L136").Value = 10319.7
$ws.Range("M136").Value = -20696.769
$ws.Range("N136").Value = -15419.7

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2187.1428
$ws.Range("I99").Value = 885
$ws.Range("J99").Value = 10000
$ws.Range("K99").Value = 885
$ws.Range("L99").Value = 10000
$ws.Range("M99").Value = 613
$ws.Range("N99").Value = -12996

$ws.Range("H139").Value = 64750
$ws.Range("J139").Value = 64750
$ws.Range("L139").Value = 64750
$ws.Range("N139").Value = -75030

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1010.55
$ws.Range("I16").Value = 1089.4
$ws.Range("J16").Value = 774
$ws.Range("K16").Value = 1089.4
$ws.Range("L16").Value = 774
$ws.Range("M16").Value = -802.4000000000001
$ws.Range("N16").Value = -1348

$ws.Range("H31").Value = 4378.4614
$ws.Range("I31").Value = 1041.1482
$ws.Range("J31").Value = 20761.637
$ws.Range("K31").Value = 1041.1482
$ws.Range("L31").Value = 20761.637
$ws.Range("M31").Value = -746.1482000000001
$ws.Range("N31").Value = -21351.637

$ws.Range("H34").Value = 4378.4614
$ws.Range("I34").Value = 1041.1482
$ws.Range("J34").Value = 20761.637
$ws.Range("K34").Value = 1041.1482
$ws.Range("L34").Value = 20761.637
$ws.Range("M34").Value = -839.1482000000001
$ws.Range("N34").Value = -21165.637

$ws.Range("H58").Value = 2571871.8
$ws.Range("I58").Value = 4360353
$ws.Range("J58").Value = 5789.7393
$ws.Range("K58").Value = 4360353
$ws.Range("L58").Value = 5789.7393
$ws.Range("M58").Value = -4360150
$ws.Range("N58").Value = -6195.7393

$ws.Range("H113").Value = 1010.55
$ws.Range("I113").Value = 1089.4
$ws.Range("J113").Value = 774
$ws.Range("K113").Value = 1089.4
$ws.Range("L113").Value = 774
$ws.Range("M113").Value = 1080.6
$ws.Range("N113").Value = -5114

$ws.Range("H132").Value = 7411324
$ws.Range("I132").Value = 15152784
$ws.Range("J132").Value = 6449.478
$ws.Range("K132").Value = 45458352
$ws.Range("L132").Value = 19348.434
$ws.Range("M132").Value = -45455822
$ws.Range("N132").Value = -24408.434

$ws.Range("H134").Value = 9192897
$ws.Range("I134").Value = 14707242
$ws.Range("J134").Value = 3678552.5
$ws.Range("K134").Value = 44121726
$ws.Range("L134").Value = 11035657.5
$ws.Range("M134").Value = -44119191
$ws.Range("N134").Value = -11040727.5

$ws.Range("H136").Value = 2571871.8
$ws.Range("I136").Value = 4360353
$ws.Range("J136").Value = 5789.7393
$ws.Range("K136").Value = 13081059
$ws.Range("L136").Value = 17369.2179
$ws.Range("M136").Value = -13078509
$ws.Range("N136").Value = -22469.2179

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 343.33334
$ws.Range("I86").Value = 378.4
$ws.Range("J86").Value = 329.84616
$ws.Range("K86").Value = 1135.2
$ws.Range("L86").Value = 989.5384799999999
$ws.Range("M86").Value = 50.80000000000018
$ws.Range("N86").Value = -3361.53848

$ws.Range("H89").Value = 343.33334
$ws.Range("I89").Value = 378.4
$ws.Range("J89").Value = 329.84616
$ws.Range("K89").Value = 3405.6
$ws.Range("L89").Value = 2968.61544
$ws.Range("M89").Value = 2522.4
$ws.Range("N89").Value = -14824.61544

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H99").Value = 11229.667
$ws.Range("I99").Value = 1585
$ws.Range("J99").Value = 30519
$ws.Range("K99").Value = 1585
$ws.Range("L99").Value = 30519
$ws.Range("M99").Value = 661
$ws.Range("N99").Value = -35011

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1657.5
$ws.Range("I68").Value = 1243.3334
$ws.Range("K68").Value = 1243.3334
$ws.Range("M68").Value = -494.3334

$ws.Range("H71").Value = 1657.5
$ws.Range("I71").Value = 1243.3334
$ws.Range("K71").Value = 6216.666999999999
$ws.Range("M71").Value = -2472.666999999999

$ws.Range("H132").Value = 4881226
$ws.Range("I132").Value = 11113838
$ws.Range("J132").Value = 3529.8696
$ws.Range("K132").Value = 33341514
$ws.Range("L132").Value = 10589.6088
$ws.Range("M132").Value = -33338984
$ws.Range("N132").Value = -15649.6088

$ws.Range("H136").Value = 12486.875
$ws.Range("I136").Value = 13982.5
$ws.Range("J136").Value = 8000
$ws.Range("K136").Value = 41947.5
$ws.Range("L136").Value = 24000
$ws.Range("M136").Value = -39397.5
$ws.Range("N136").Value = -29100

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 23066474
$ws.Range("I132").Value = 11827601
$ws.Range("J132").Value = 36713676
$ws.Range("K132").Value = 35482803
$ws.Range("L132").Value = 110141028
$ws.Range("M132").Value = -35480273
$ws.Range("N132").Value = -110146088

$ws.Range("H136").Value = 20407976
$ws.Range("I136").Value = 13735093
$ws.Range("J136").Value = 31251412
$ws.Range("K136").Value = 41205279
$ws.Range("L136").Value = 93754236
$ws.Range("M136").Value = -41202729
$ws.Range("N136").Value = -93759336
